{"js": "// Soften small-space food ebook tone to framework (v0.3)\n// Each entry is [findText, replaceText]; both sides are matched\n// case-sensitively against exact paragraph text captured in the diff.\nconst replacements = [\n  [\"First: **it\u2019s physics, not vibes.** Containers are small soil volumes exposed to heat and wind. Miss watering during a hot spell and production drops\u2014often for the rest of that fruiting cycle.\", \"First: **it\u2019s physics, not vibes.** Containers are small soil volumes exposed to heat and wind. A short watering failure during a hot spell can depress production for the rest of that fruiting cycle.\"],\n  [\"Third: **logistics beat technique.** Soil volume, irrigation reliability, pest exclusion, and replacement planting determine yield. Beginners overspend on seeds and underspend on containers, mix quality, and watering infrastructure.\", \"Third: **logistics beat technique.** Soil volume, irrigation reliability, pest exclusion, and replacement planting determine yield. Beginners overspend on seeds and underspend on containers, mix quality, and watering reliability.\"],\n  [\"**Define \u201cmeaningful\u201d upfront:** weigh harvest weekly; if you don\u2019t weigh, you don\u2019t know.\", \"**Define \u201cmeaningful\u201d upfront:** weekly harvest weight is the only honest scoreboard.\"],\n  [\"**Solve water before crops:** inconsistent water guarantees inconsistent yield.\", \"**Water reliability sits upstream of everything:** when water is inconsistent, everything else becomes noise.\"],\n  [\"**Use fewer, larger containers:** 15\u201325 gallons beats many small pots; small pots are drought machines.\", \"**Fewer, larger containers change the game:** 15\u201325 gallons outperforms many small pots; small pots behave like drought simulators.\"],\n  [\"**Run a \u201ccalories + nutrition\u201d mix:** at least one calorie crop (tubers/beans/squash) plus one nutrient green.\", \"**A \u201ccalories + nutrition\u201d mix is the only mix that feels meaningful:** at least one calorie crop (tubers/beans/squash) plus one nutrient green.\"],\n  [\"**Treat sun hours as a gate:** <6 hours direct sun means greens/beans unless you add grow lights.\", \"**Sun hours are a constraint, not a suggestion:** <6 hours direct sun shifts you toward greens/beans unless you add grow lights.\"],\n  [\"**Install pest exclusion early:** netting/row cover prevents the late-response spiral.\", \"**Pest exclusion is the adult move:** netting/row cover prevents the late-response spiral.\"],\n  [\"**Trellis for access and plant health:** vertical structure is control, not magic.\", \"**Trellising is control:** it protects access, airflow, and harvesting more than it boosts yield.\"],\n  [\"**Replace weak plants quickly:** stalled for 2+ weeks in prime season = replace.\", \"**Replacement planting is normal:** a plant that stalls for weeks is an output risk, not a project.\"],\n  [\"**Stagger planting for continuity:** 2\u20133 waves of greens/beans beats one peak-and-crash planting.\", \"**Staggering keeps output continuous:** 2\u20133 waves of greens/beans beats one peak-and-crash planting.\"],\n  [\"**Mulch containers aggressively:** 2\u20134 inches to reduce temperature swings and watering load.\", \"**Mulch buys stability:** 2\u20134 inches reduces temperature swings and watering load.\"],\n  [\"**Feed lightly and regularly:** avoid \u201chero feeding\u201d after neglect; it creates problems.\", \"**Light, regular feeding beats rescue feeding:** neglect + \u201chero fertilizer\u201d is a reliable way to invite pests and disorder.\"],\n  [\"Decision: <6 hours direct sun and no grow lights \u2192 plan for supplemental vegetables, not calorie coverage.\", \"What matters: if you don\u2019t have \u22656 hours direct sun (or equivalent light), \u201cmeaningful\u201d will come from greens/beans, not calorie crops.\"],\n  [\"Decision: 15\u201325 gallons for heavy fruiters; 7\u201310 gallons is a compromise.\", \"What matters: heavy fruiters need real volume; 15\u201325 gallons is where they stop acting fragile.\"],\n  [\"Decision: if you miss watering >1\u00d7/month in summer, install drip + timer or cut plant count.\", \"What matters: if you routinely miss waterings in summer, scale the system to your life or automate; willpower doesn\u2019t fix heat.\"],\n  [\"Decision: pick crops that match your heat/cold profile; stop fighting your climate.\", \"What matters: you want to learn *one system*, not collect experiences.\"],\n  [\"Decision: potatoes/sweet potatoes and dry beans usually beat tomatoes on calories per effort.\", \"What matters: potatoes/sweet potatoes and dry beans usually win on calories per effort; tomatoes win on satisfaction, not calories.\"],\n  [\"Decision: if pests are common where you live, start with exclusion and escalate only when needed.\", \"What matters: if pests are common where you live, \u201creactive\u201d becomes a seasonal tax; exclusion is how you keep momentum.\"],\n  [\"1) **Set your output target and measurement**\", \"1) **Define what \u201cmeaningful\u201d means in your household**\"],\n  [\"Intent: turn \u201cgardening\u201d into supply. Weekly weigh-ins and a one-page log.\", \"What matters at this stage is measurement. If weight isn\u2019t tracked weekly, the system drifts into vibes.\"],\n  [\"2) **Lock constraints (space, light, travel, budget)**\", \"2) **Make constraints explicit**\"],\n  [\"Intent: prevent plans that fail on your calendar. Travel requires automation or fewer plants.\", \"What matters at this stage is honesty: light, time, travel, and budget. A plan that conflicts with your calendar will fail quietly.\"],\n  [\"3) **Build the foundation (containers/bed + soil + mulch + water plan)**\", \"3) **Build stability before ambition**\"],\n  [\"Intent: remove failure points before planting. Don\u2019t start until water and soil are stable.\", \"What matters at this stage is removing predictable failure points (soil volume, mulch, water reliability). Seeds are easy; stability is the work.\"],\n  [\"4) **Plant a focused crop set with sequencing**\", \"4) **Run a focused crop set with continuity**\"],\n  [\"Intent: continuity. Calorie crops early; greens/beans in waves; reserve space for replacements.\", \"What matters at this stage is consistent output: staggered planting and space reserved for replacements.\"],\n  [\"5) **Operate like maintenance**\", \"5) **Operate like maintenance, not emergencies**\"],\n  [\"Intent: avoid rescue cycles. Standardize watering, scouting, harvesting, and replanting.\", \"What matters at this stage is rhythm: water, scout, harvest, replant. Rescue cycles consume time and produce less.\"],\n  [\"6) **Cull and reallocate by evidence**\", \"6) **Reallocate by evidence**\"],\n  [\"Intent: protect yield. Underperformers lose their spot; winners get more space.\", \"What matters at this stage is detachment: underperformers lose space; reliable performers earn space.\"],\n  [\"Don\u2019t assume \u201cmore fertilizer\u201d fixes stress.\", \"Don\u2019t treat fertilizer as a substitute for stable water.\"],\n  [\"Don\u2019t keep weak plants \u201cto see if they recover.\u201d Replace them.\", \"Don\u2019t keep weak plants \u201cto see if they recover.\u201d\"],\n  [\"Pests twice in one week = switch to exclusion.\", \"Pests twice in one week = exclusion time.\"],\n  [\"Missing waterings weekly \u2192 reduce plant count until you hit 0 missed waterings/month.\", \"Missing waterings weekly \u2192 shrink the system until misses stop.\"],\n  [\"Repeated wilting in one location \u2192 move containers or add shade/wind control.\", \"Repeated wilting in one location \u2192 treat the location as the problem: move containers or add shade/wind control.\"],\n  [\"Escalating pest interventions with continued losses \u2192 switch to exclusion or change crops.\", \"Escalating pest interventions with continued losses \u2192 change posture (exclusion) or change crops.\"],\n  [\"Lots of leaves, little food weight \u2192 reallocate to calorie crops or larger containers.\", \"Lots of leaves, little food weight \u2192 admit the mismatch: reallocate to calorie crops or larger containers.\"],\n  [\"Inputs (soil/water/time) remain unstable \u2192 stop scaling and run a smaller stable system.\", \"Inputs (soil/water/time) remain unstable \u2192 stop scaling; run a smaller stable system until it\u2019s boring.\"],\n  [\"Small-space food is operations: measurement, reliability, replacement, sequencing. Run it that way and the harvest becomes predictable.\", \"Small-space food is operations: measurement, reliability, replacement, sequencing. Structure orients you; consistency delivers the harvest.\"],\n];\n\nconst body = context.document.body;\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + findText);\n  }\n  for (const result of results.items) {\n    result.insertText(replaceText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Soften small-space food ebook tone to framework (v0.3)\n$d = $word.ActiveDocument\n\nfunction Replace-Text($FindText, $ReplaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($FindText, $false, $true, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $FindText\"\n    }\n}\n\n$replacements = @(\n    @(\"First: **it\u2019s physics, not vibes.** Containers are small soil volumes exposed to heat and wind. Miss watering during a hot spell and production drops\u2014often for the rest of that fruiting cycle.\", \"First: **it\u2019s physics, not vibes.** Containers are small soil volumes exposed to heat and wind. A short watering failure during a hot spell can depress production for the rest of that fruiting cycle.\"),\n    @(\"Third: **logistics beat technique.** Soil volume, irrigation reliability, pest exclusion, and replacement planting determine yield. Beginners overspend on seeds and underspend on containers, mix quality, and watering infrastructure.\", \"Third: **logistics beat technique.** Soil volume, irrigation reliability, pest exclusion, and replacement planting determine yield. Beginners overspend on seeds and underspend on containers, mix quality, and watering reliability.\"),\n    @(\"**Define \u201cmeaningful\u201d upfront:** weigh harvest weekly; if you don\u2019t weigh, you don\u2019t know.\", \"**Define \u201cmeaningful\u201d upfront:** weekly harvest weight is the only honest scoreboard.\"),\n    @(\"**Solve water before crops:** inconsistent water guarantees inconsistent yield.\", \"**Water reliability sits upstream of everything:** when water is inconsistent, everything else becomes noise.\"),\n    @(\"**Use fewer, larger containers:** 15\u201325 gallons beats many small pots; small pots are drought machines.\", \"**Fewer, larger containers change the game:** 15\u201325 gallons outperforms many small pots; small pots behave like drought simulators.\"),\n    @(\"**Run a \u201ccalories + nutrition\u201d mix:** at least one calorie crop (tubers/beans/squash) plus one nutrient green.\", \"**A \u201ccalories + nutrition\u201d mix is the only mix that feels meaningful:** at least one calorie crop (tubers/beans/squash) plus one nutrient green.\"),\n    @(\"**Treat sun hours as a gate:** <6 hours direct sun means greens/beans unless you add grow lights.\", \"**Sun hours are a constraint, not a suggestion:** <6 hours direct sun shifts you toward greens/beans unless you add grow lights.\"),\n    @(\"**Install pest exclusion early:** netting/row cover prevents the late-response spiral.\", \"**Pest exclusion is the adult move:** netting/row cover prevents the late-response spiral.\"),\n    @(\"**Trellis for access and plant health:** vertical structure is control, not magic.\", \"**Trellising is control:** it protects access, airflow, and harvesting more than it boosts yield.\"),\n    @(\"**Replace weak plants quickly:** stalled for 2+ weeks in prime season = replace.\", \"**Replacement planting is normal:** a plant that stalls for weeks is an output risk, not a project.\"),\n    @(\"**Stagger planting for continuity:** 2\u20133 waves of greens/beans beats one peak-and-crash planting.\", \"**Staggering keeps output continuous:** 2\u20133 waves of greens/beans beats one peak-and-crash planting.\"),\n    @(\"**Mulch containers aggressively:** 2\u20134 inches to reduce temperature swings and watering load.\", \"**Mulch buys stability:** 2\u20134 inches reduces temperature swings and watering load.\"),\n    @(\"**Feed lightly and regularly:** avoid \u201chero feeding\u201d after neglect; it creates problems.\", \"**Light, regular feeding beats rescue feeding:** neglect + \u201chero fertilizer\u201d is a reliable way to invite pests and disorder.\"),\n    @(\"Decision: <6 hours direct sun and no grow lights \u2192 plan for supplemental vegetables, not calorie coverage.\", \"What matters: if you don\u2019t have \u22656 hours direct sun (or equivalent light), \u201cmeaningful\u201d will come from greens/beans, not calorie crops.\"),\n    @(\"Decision: 15\u201325 gallons for heavy fruiters; 7\u201310 gallons is a compromise.\", \"What matters: heavy fruiters need real volume; 15\u201325 gallons is where they stop acting fragile.\"),\n    @(\"Decision: if you miss watering >1\u00d7/month in summer, install drip + timer or cut plant count.\", \"What matters: if you routinely miss waterings in summer, scale the system to your life or automate; willpower doesn\u2019t fix heat.\"),\n    @(\"Decision: pick crops that match your heat/cold profile; stop fighting your climate.\", \"What matters: you want to learn *one system*, not collect experiences.\"),\n    @(\"Decision: potatoes/sweet potatoes and dry beans usually beat tomatoes on calories per effort.\", \"What matters: potatoes/sweet potatoes and dry beans usually win on calories per effort; tomatoes win on satisfaction, not calories.\"),\n    @(\"Decision: if pests are common where you live, start with exclusion and escalate only when needed.\", \"What matters: if pests are common where you live, \u201creactive\u201d becomes a seasonal tax; exclusion is how you keep momentum.\"),\n    @(\"1) **Set your output target and measurement**\", \"1) **Define what \u201cmeaningful\u201d means in your household**\"),\n    @(\"Intent: turn \u201cgardening\u201d into supply. Weekly weigh-ins and a one-page log.\", \"What matters at this stage is measurement. If weight isn\u2019t tracked weekly, the system drifts into vibes.\"),\n    @(\"2) **Lock constraints (space, light, travel, budget)**\", \"2) **Make constraints explicit**\"),\n    @(\"Intent: prevent plans that fail on your calendar. Travel requires automation or fewer plants.\", \"What matters at this stage is honesty: light, time, travel, and budget. A plan that conflicts with your calendar will fail quietly.\"),\n    @(\"3) **Build the foundation (containers/bed + soil + mulch + water plan)**\", \"3) **Build stability before ambition**\"),\n    @(\"Intent: remove failure points before planting. Don\u2019t start until water and soil are stable.\", \"What matters at this stage is removing predictable failure points (soil volume, mulch, water reliability). Seeds are easy; stability is the work.\"),\n    @(\"4) **Plant a focused crop set with sequencing**\", \"4) **Run a focused crop set with continuity**\"),\n    @(\"Intent: continuity. Calorie crops early; greens/beans in waves; reserve space for replacements.\", \"What matters at this stage is consistent output: staggered planting and space reserved for replacements.\"),\n    @(\"5) **Operate like maintenance**\", \"5) **Operate like maintenance, not emergencies**\"),\n    @(\"Intent: avoid rescue cycles. Standardize watering, scouting, harvesting, and replanting.\", \"What matters at this stage is rhythm: water, scout, harvest, replant. Rescue cycles consume time and produce less.\"),\n    @(\"6) **Cull and reallocate by evidence**\", \"6) **Reallocate by evidence**\"),\n    @(\"Intent: protect yield. Underperformers lose their spot; winners get more space.\", \"What matters at this stage is detachment: underperformers lose space; reliable performers earn space.\"),\n    @(\"Don\u2019t assume \u201cmore fertilizer\u201d fixes stress.\", \"Don\u2019t treat fertilizer as a substitute for stable water.\"),\n    @(\"Don\u2019t keep weak plants \u201cto see if they recover.\u201d Replace them.\", \"Don\u2019t keep weak plants \u201cto see if they recover.\u201d\"),\n    @(\"Pests twice in one week = switch to exclusion.\", \"Pests twice in one week = exclusion time.\"),\n    @(\"Missing waterings weekly \u2192 reduce plant count until you hit 0 missed waterings/month.\", \"Missing waterings weekly \u2192 shrink the system until misses stop.\"),\n    @(\"Repeated wilting in one location \u2192 move containers or add shade/wind control.\", \"Repeated wilting in one location \u2192 treat the location as the problem: move containers or add shade/wind control.\"),\n    @(\"Escalating pest interventions with continued losses \u2192 switch to exclusion or change crops.\", \"Escalating pest interventions with continued losses \u2192 change posture (exclusion) or change crops.\"),\n    @(\"Lots of leaves, little food weight \u2192 reallocate to calorie crops or larger containers.\", \"Lots of leaves, little food weight \u2192 admit the mismatch: reallocate to calorie crops or larger containers.\"),\n    @(\"Inputs (soil/water/time) remain unstable \u2192 stop scaling and run a smaller stable system.\", \"Inputs (soil/water/time) remain unstable \u2192 stop scaling; run a smaller stable system until it\u2019s boring.\"),\n    @(\"Small-space food is operations: measurement, reliability, replacement, sequencing. Run it that way and the harvest becomes predictable.\", \"Small-space food is operations: measurement, reliability, replacement, sequencing. Structure orients you; consistency delivers the harvest.\"),\n)\n\nforeach ($pair in $replacements) {\n    Replace-Text $pair[0] $pair[1]\n}\n\nWrite-Output \"Applied $($replacements.Count) replacements\"\n"}
